$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "42.849.54"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -1.25%  "
$c.ClearFormats()
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.325.81"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +0.83%  "
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.ClearFormats()
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "304.73"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -1.97%  "
$c.ClearFormats()
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "100.55"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -2.97%  "
$c.ClearFormats()
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.508"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -4.71%  "
$c.ClearFormats()
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.505"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -4.34%  "
$c.ClearFormats()
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.51"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -5.96%  "
$c.ClearFormats()
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "52.20"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -0.39%  "
$c.ClearFormats()
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0794"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -2.17%  "
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +0.95%  "
$c.ClearFormats()
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.76"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -3.77%  "
$c.ClearFormats()
$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = "WrappedliquidstakedEther2.0"
$c.ClearFormats()
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.ClearFormats()
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.686.76"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +0.91%  "
$c.ClearFormats()
$c = $ws.Range("B16")
$c.NumberFormat = "@"
$c.Value = "Chainlink"
$c.ClearFormats()
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c.ClearFormats()
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "15.70"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +4.07%  "
$c.ClearFormats()
$c = $ws.Range("B17")
$c.NumberFormat = "@"
$c.Value = "WrappedEther"
$c.ClearFormats()
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c.ClearFormats()
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.342.84"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +1.12%  "
$c.ClearFormats()
$c = $ws.Range("B18")
$c.NumberFormat = "@"
$c.Value = "Polygon"
$c.ClearFormats()
$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c.ClearFormats()
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.821"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +1.39%  "
$c.ClearFormats()
$c = $ws.Range("B19")
$c.NumberFormat = "@"
$c.Value = "WrappedBTC"
$c.ClearFormats()
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c.ClearFormats()
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "42.743.03"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -1.24%  "
$c.ClearFormats()
$c = $ws.Range("B20")
$c.NumberFormat = "@"
$c.Value = "ShibaInu"
$c.ClearFormats()
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c.ClearFormats()
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0904"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  -2.41%  "
$c.ClearFormats()
$c = $ws.Range("B21")
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.ClearFormats()
$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.ClearFormats()
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -0.41%  "
$c.ClearFormats()
$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = "InternetComputer(DFINITY)"
$c.ClearFormats()
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c.ClearFormats()
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "11.61"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -4.62%  "
$c.ClearFormats()
$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = "Litecoin"
$c.ClearFormats()
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c.ClearFormats()
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "69.30"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +1.76%  "
$c.ClearFormats()
$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.ClearFormats()
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "235.61"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -2.90%  "
$c.ClearFormats()
$c = $ws.Range("B25")
$c.NumberFormat = "@"
$c.Value = "ImmutableX"
$c.ClearFormats()
$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c.ClearFormats()
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -1.73%  "
$c.ClearFormats()
$c = $ws.Range("B26")
$c.NumberFormat = "@"
$c.Value = "PancakeSwap"
$c.ClearFormats()
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c.ClearFormats()
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  -3.37%  "
$c.ClearFormats()
$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = "Dai"
$c.ClearFormats()
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.08%  "
$c.ClearFormats()
$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = "EthereumClassic"
$c.ClearFormats()
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "25.41"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +1.97%  "
$c.ClearFormats()
$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = "LEO"
$c.ClearFormats()
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c.ClearFormats()
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "3.96"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -0.56%  "
$c.ClearFormats()
$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = "Toncoin"
$c.ClearFormats()
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c.ClearFormats()
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -1.33%  "
$c.ClearFormats()
$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = "InjectiveProtocol"
$c.ClearFormats()
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "34.77"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -5.75%  "
$c.ClearFormats()
$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = "Cosmos"
$c.ClearFormats()
$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "9.22"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -4.46%  "
$c.ClearFormats()
$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = "Monero"
$c.ClearFormats()
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "160.24"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -4.51%  "
$c.ClearFormats()
$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = "FirstDigitalUSD"
$c.ClearFormats()
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.ClearFormats()
$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = "Filecoin"
$c.ClearFormats()
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c.ClearFormats()
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.06"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -4.14%  "
$c.ClearFormats()
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +3.40%  "
$c.ClearFormats()
$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "WEMIXToken"
$c.ClearFormats()
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c.ClearFormats()
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -3.38%  "
$c.ClearFormats()
$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "Hedera"
$c.ClearFormats()
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.ClearFormats()
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0720"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -3.09%  "
$c.ClearFormats()
$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = "Celestia"
$c.ClearFormats()
$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c.ClearFormats()
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "17.00"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -6.50%  "
$c.ClearFormats()
$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = "LidoDAOToken"
$c.ClearFormats()
$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -5.53%  "
$c.ClearFormats()
$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = "ARBITRUM"
$c.ClearFormats()
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c.ClearFormats()
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -2.74%  "
$c.ClearFormats()
$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = "Kaspa"
$c.ClearFormats()
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c.ClearFormats()
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.101"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -4.67%  "
$c.ClearFormats()
$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c.ClearFormats()
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.ClearFormats()
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -3.20%  "
$c.ClearFormats()
$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = "ApeXProtocol"
$c.ClearFormats()
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c.ClearFormats()
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.55"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -6.08%  "
$c.ClearFormats()
$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = "Maker"
$c.ClearFormats()
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c.ClearFormats()
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.004.88"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +1.31%  "
$c.ClearFormats()
$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = "VeChain"
$c.ClearFormats()
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c.ClearFormats()
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0282"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -4.05%  "
$c.ClearFormats()
$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = "EnergySwap"
$c.ClearFormats()
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.ClearFormats()
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "18.75"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  -1.59%  "
$c.ClearFormats()
$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = "FraxShare"
$c.ClearFormats()
$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.ClearFormats()
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "10.16"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.81%  "
$c.ClearFormats()
$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c.ClearFormats()
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -4.10%  "
$c.ClearFormats()
$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = "MultiversX"
$c.ClearFormats()
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c.ClearFormats()
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "55.44"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -0.87%  "
$c.ClearFormats()
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = "HuobiToken"
$c.ClearFormats()
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.86"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -2.23%  "
$c.ClearFormats()
